# "update all except prostate"
# The Prostate cohort rows (v1.1-consortium, 9 dataset rows with no
# synapse_id / release_date yet) are removed entirely; every row below
# them shifts up to fill the gap, so the sheet shrinks from A1:E115 to
# A1:E106.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 89-97 hold the nine "Prostate" / "v1.1-consortium" placeholder
# rows (cancer_level_dataset_index ... tm_level_dataset) that have no
# synapse_id or release_date filled in yet. Delete them outright.
$ws.Range("A89:E97").EntireRow.Delete()

# Reflect the editor's final cursor position/selection from the diff
# (topLeftCell scrolled to A82, active cell C97 after the shift).
$ws.Application.Goto($ws.Range("A82"))
$ws.Range("C97").Select()
